$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: update the date in A1 (45308 -> 45309, i.e. +1 day)
$ws.Range("A1").Value = 45309

# Step 2: update the prices in D30 and D31
$ws.Range("D30").Value = 760
$ws.Range("D31").Value = 1520
